$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# Preserve the sheet's pre-existing (stale) dimension extent through column AB
# (the original workbook's <dimension> already spanned to AB even though AB
# itself holds no data; touching it keeps that extent after our edits).
$ws.Range("AB1").Font.Bold = $false

# "nexial.web.dragFrom" / dragTo(...) is a new "web" command; the #system sheet
# keeps an alphabetically-sorted catalogue of web commands in column V, which
# backs the "web" named range. Insert the new entry and re-flow everything
# that sorts after it by one row (V63:V120).
$wb.Names.Item("web").RefersTo = "='#system'!`$V`$2:`$V`$120"
$ws.Range("V63").Value = "dragTo(fromLocator,xOffset,yOffset)"
$ws.Range("V64").Value = "editLocalStorage(key,value)"
$ws.Range("V65").Value = "executeScript(var,script)"
$ws.Range("V66").Value = "focus(locator)"
$ws.Range("V67").Value = "goBack()"
$ws.Range("V68").Value = "goBackAndWait()"
$ws.Range("V69").Value = "maximizeWindow()"
$ws.Range("V70").Value = "mouseOver(locator)"
$ws.Range("V71").Value = "open(url)"
$ws.Range("V72").Value = "openAndWait(url,waitMs)"
$ws.Range("V73").Value = "openHttpBasic(url,username,password)"
$ws.Range("V74").Value = "openIgnoreTimeout(url)"
$ws.Range("V75").Value = "refresh()"
$ws.Range("V76").Value = "refreshAndWait()"
$ws.Range("V77").Value = "resizeWindow(width,height)"
$ws.Range("V78").Value = "saveAllWindowIds(var)"
$ws.Range("V79").Value = "saveAllWindowNames(var)"
$ws.Range("V80").Value = "saveAttribute(var,locator,attrName)"
$ws.Range("V81").Value = "saveCount(var,locator)"
$ws.Range("V82").Value = "saveDivsAsCsv(headers,rows,cells,nextPage,file)"
$ws.Range("V83").Value = "saveElement(var,locator)"
$ws.Range("V84").Value = "saveElements(var,locator)"
$ws.Range("V85").Value = "saveLocalStorage(var,key)"
$ws.Range("V86").Value = "saveLocation(var)"
$ws.Range("V87").Value = "savePageAs(var,sessionIdName,url)"
$ws.Range("V88").Value = "savePageAsFile(sessionIdName,url,file)"
$ws.Range("V89").Value = "saveTableAsCsv(locator,nextPageLocator,file)"
$ws.Range("V90").Value = "saveText(var,locator)"
$ws.Range("V91").Value = "saveTextArray(var,locator)"
$ws.Range("V92").Value = "saveTextSubstringAfter(var,locator,delim)"
$ws.Range("V93").Value = "saveTextSubstringBefore(var,locator,delim)"
$ws.Range("V94").Value = "saveTextSubstringBetween(var,locator,start,end)"
$ws.Range("V95").Value = "saveValue(var,locator)"
$ws.Range("V96").Value = "scrollLeft(locator,pixel)"
$ws.Range("V97").Value = "scrollRight(locator,pixel)"
$ws.Range("V98").Value = "scrollTo(locator)"
$ws.Range("V99").Value = "select(locator,text)"
$ws.Range("V100").Value = "selectFrame(locator)"
$ws.Range("V101").Value = "selectMulti(locator,array)"
$ws.Range("V102").Value = "selectMultiOptions(locator)"
$ws.Range("V103").Value = "selectText(locator)"
$ws.Range("V104").Value = "selectWindow(winId)"
$ws.Range("V105").Value = "selectWindowAndWait(winId,waitMs)"
$ws.Range("V106").Value = "selectWindowByIndex(index)"
$ws.Range("V107").Value = "selectWindowByIndexAndWait(index,waitMs)"
$ws.Range("V108").Value = "toggleSelections(locator)"
$ws.Range("V109").Value = "type(locator,value)"
$ws.Range("V110").Value = "typeKeys(locator,value)"
$ws.Range("V111").Value = "uncheckAll(locator)"
$ws.Range("V112").Value = "unselectAllText()"
$ws.Range("V113").Value = "upload(fieldLocator,file)"
$ws.Range("V114").Value = "verifyContainText(locator,text)"
$ws.Range("V115").Value = "verifyText(locator,text)"
$ws.Range("V116").Value = "wait(waitMs)"
$ws.Range("V117").Value = "waitForElementPresent(locator)"
$ws.Range("V118").Value = "waitForPopUp(winId,waitMs)"
$ws.Range("V119").Value = "waitForTextPresent(text)"
$ws.Range("V120").Value = "waitForTitle(text)"

# Two new "xml" commands (beautify/minify) are catalogued in column AA, which
# backs the "xml" named range. Insert them alphabetically and re-flow the
# three existing store* entries down by two rows (AA9:AA13).
$wb.Names.Item("xml").RefersTo = "='#system'!`$AA`$2:`$AA`$13"
$ws.Range("AA9").Value = "beautify(xml,var)"
$ws.Range("AA10").Value = "minify(xml,var)"
$ws.Range("AA11").Value = "storeCount(xml,xpath,var)"
$ws.Range("AA12").Value = "storeValue(xml,xpath,var)"
$ws.Range("AA13").Value = "storeValues(xml,xpath,var)"
